# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1132038
$ws.Range("C4").Value = 1008
$ws.Range("E4").Value = 904589
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 65783

# --- Row 20: Suiza ---
$ws.Range("F20").Value = 141

# --- Row 24: Suecia ---
$ws.Range("B24").Value = 22082
$ws.Range("C24").Value = 562
$ws.Range("E24").Value = 18408
$ws.Range("G24").Value = 16
$ws.Range("H24").Value = 2669

# --- Row 72: Croacia ---
$ws.Range("B72").Value = 2088
$ws.Range("C72").Value = 3
$ws.Range("D72").Value = 1463
$ws.Range("E72").Value = 548
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 77

# --- Insert new country "Bosnia y Herzegovina" at row 75, pushing ---
# --- Camerun and Islandia down by one row (shared-string table gets ---
# --- reordered so Bosnia y Herzegovina sits right after Azerbaiyan) ---
$ws.Range("A75").Value = "Bosnia y Herzegovina"
$ws.Range("B75").Value = 1839
$ws.Range("C75").Value = 58
$ws.Range("D75").Value = 779
$ws.Range("E75").Value = 988
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 72

$ws.Range("A76").Value = "Camerun"
$ws.Range("B76").Value = 1832
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 934
$ws.Range("E76").Value = 837
$ws.Range("F76").Value = 12
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 61

$ws.Range("A77").Value = "Islandia"
$ws.Range("B77").Value = 1798
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 1689
$ws.Range("E77").Value = 99
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 10

# --- Row 128: Isla de Man ---
$ws.Range("B128").Value = 320
$ws.Range("C128").Value = 4
$ws.Range("E128").Value = 27

# --- Row 186: Laos ---
$ws.Range("D186").Value = 9
$ws.Range("E186").Value = 10
